$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# A new data row (row 9) was appended below the existing "Etat Taxes" table
# (rows 1-8). Only columns A-H are populated for this entry; I-M are left
# blank, same as in the committed workbook.
$ws.Range("A9").Value = "001/LF/DR IFRAN"
$ws.Range("B9").Value = "Logement de fonction"
$ws.Range("C9").Value = "mdkjhf"
$ws.Range("D9").Value = "test test"
$ws.Range("E9").Value = "ds"
$ws.Range("F9").Value = "mensuelle"
$ws.Range("G9").Value = 0
$ws.Range("H9").Value = 0

# The sheet ignores "number stored as text" warnings across A1:M8; extend
# that same ignore-error range to cover the newly added row 9 as well
# (xlNumberAsText = 3), matching Excel's own "Ignore Error" behaviour when
# a user edits an existing ignored range.
$ws.Range("A1:M9").Errors.Item(3).Ignore = $true
